$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = -13.53519999999999
$ws.Range("C4").Value = -14.38140000000001
$ws.Range("D6").Value = -7.872300000000001
$ws.Range("C7").Value = -11.36749999999999
$ws.Range("D7").Value = -7.599399999999997
$ws.Range("C8").Value = -12.14789999999998
$ws.Range("D8").Value = -8.70619999999999
$ws.Range("A11").Value = -22.02160000000002
$ws.Range("E11").Value = 13.65380000000001
$ws.Range("A12").Value = -21.04100000000001
$ws.Range("C12").Value = -11.28839999999999
$ws.Range("C14").Value = -12.09250000000001
$ws.Range("E14").Value = 14.0155
$ws.Range("A15").Value = -21.18010000000002
$ws.Range("D19").Value = -7.613299999999997
$ws.Range("E19").Value = 14.08060000000001
$ws.Range("D21").Value = -8.621599999999994
$ws.Range("E21").Value = 12.54070000000001
$ws.Range("C22").Value = -11.06159999999999
$ws.Range("D24").Value = -7.727999999999987
$ws.Range("D25").Value = -7.527899999999999
